$d = $word.ActiveDocument

# Position at the very end of the document (end of the last paragraph,
# "Made a new projects.html page ...").
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)

# --- 1) New empty paragraph ------------------------------------------------
# Mirrors the blank separator paragraphs used elsewhere between diary
# entries. A freshly-inserted empty paragraph carries a placeholder run,
# so type a throw-away character into it and delete just that character
# (not the paragraph mark) to end up with a bare <w:p> - no run at all -
# exactly like the existing blank paragraphs in this document.
$r.InsertParagraphAfter()
$rEmpty = $d.Paragraphs.Last.Range
$rEmpty.Collapse(0)
$rEmpty.InsertAfter("X")
$pEmpty = $d.Paragraphs.Last
$trimRange = $d.Range($pEmpty.Range.Start, $pEmpty.Range.Start + 1)
$trimRange.Delete()

# --- 2) Date paragraph "28.04.2021" ----------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertAfter("28.04.2021")

# --- 3) Entry paragraph made of two runs -----------------------------------
# Both runs share identical formatting, but the original document keeps
# consecutive typing sessions as separate <w:r> elements. Build each half
# as its own paragraph, then remove the paragraph mark between them to
# splice them into one paragraph containing two runs.
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pA = $d.Paragraphs.Last
$rA = $pA.Range
$rA.Collapse(0)
$rA.InsertAfter("Making my final finishing touches. I am going to add some javascript to this site and utilize a modal. Found an example from codepen.io and see if that fits with a bit of ")

$rA2 = $d.Paragraphs.Last.Range
$rA2.Collapse(0)
$rA2.InsertParagraphAfter()
$pB = $d.Paragraphs.Last
$rB = $pB.Range
$rB.Collapse(0)
$rB.InsertAfter("modifying.")

$joinStart = $pA.Range.End - 1
$joinRange = $d.Range($joinStart, $pB.Range.Start)
$joinRange.Delete()

# --- 4) Final paragraph of the entry ---------------------------------------
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertAfter("I am finding problematic to open all 3 buttons with same code and names. I also saw that Materialize has its own modals. Going to test those ones as well.")

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
